# Applies the scheduled-runner value refresh to the Lamia_Profits workbook.
# For each affected row, columns H/I/J/K/L/M/N are recomputed pricing-analysis
# figures; some rows gain or lose their M/N cell entirely (Excel omits cells that
# hold no value, so ClearContents() is used where a cell must disappear).
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 507.5
$ws.Range("I33").Value = 437.42856
$ws.Range("K33").Value = 437.42856
$ws.Range("M33").Value = -208.42856
$ws.Range("H86").Value = 6231.231
$ws.Range("I86").Value = 6667.3335
$ws.Range("J86").Value = 6100.4
$ws.Range("K86").Value = 6667.3335
$ws.Range("L86").Value = 6100.4
$ws.Range("M86").Value = -5544.3335
$ws.Range("N86").Value = -8346.4
$ws.Range("H89").Value = 6231.231
$ws.Range("I89").Value = 6667.3335
$ws.Range("J89").Value = 6100.4
$ws.Range("K89").Value = 33336.6675
$ws.Range("L89").Value = 30502
$ws.Range("M89").Value = -27720.6675
$ws.Range("N89").Value = -41734
$ws.Range("H112").Value = 1328.16
$ws.Range("J112").Value = 1360.1818
$ws.Range("L112").Value = 4080.5454
$ws.Range("N112").Value = -6296.5454
$ws.Range("H132").Value = 969
$ws.Range("I132").Value = 741.1613
$ws.Range("K132").Value = 2223.4839
$ws.Range("M132").Value = 306.5160999999998
$ws.Range("H138").Value = 2961.5732
$ws.Range("I138").Value = 1694.1364
$ws.Range("K138").Value = 5082.4092
$ws.Range("M138").Value = 57.59079999999994

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7516.533
$ws.Range("I2").Value = 761.5
$ws.Range("J2").Value = 34536.668
$ws.Range("K2").Value = 761.5
$ws.Range("L2").Value = 34536.668
$ws.Range("M2").Value = -648.5
$ws.Range("N2").Value = -34762.668
$ws.Range("H19").Value = 19999
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("H26").Value = 2058.1667
$ws.Range("I26").Value = 2483
$ws.Range("K26").Value = 2483
$ws.Range("M26").Value = -2153
$ws.Range("H36").Value = 15713.857
$ws.Range("I36").Value = 9999
$ws.Range("K36").Value = 9999
$ws.Range("M36").Value = -9653
$ws.Range("H45").Value = 58825560
$ws.Range("I45").Value = 62502100
$ws.Range("K45").Value = 62502100
$ws.Range("M45").Value = -62501723
$ws.Range("H61").Value = 7129.2964
$ws.Range("I61").Value = 6342.1055
$ws.Range("J61").Value = 8998.875
$ws.Range("K61").Value = 6342.1055
$ws.Range("L61").Value = 8998.875
$ws.Range("M61").Value = -6130.1055
$ws.Range("N61").Value = -9422.875
$ws.Range("H74").Value = 15153861
$ws.Range("I74").Value = 15874568
$ws.Range("K74").Value = 15874568
$ws.Range("M74").Value = -15873694
$ws.Range("H77").Value = 15153861
$ws.Range("I77").Value = 15874568
$ws.Range("K77").Value = 79372840
$ws.Range("M77").Value = -79368472
$ws.Range("H97").Value = 1507.7778
$ws.Range("I97").Value = 1507.7778
$ws.Range("K97").Value = 1507.7778
$ws.Range("M97").Value = -1011.7778
$ws.Range("H116").Value = 7516.533
$ws.Range("I116").Value = 761.5
$ws.Range("J116").Value = 34536.668
$ws.Range("K116").Value = 761.5
$ws.Range("L116").Value = 34536.668
$ws.Range("M116").Value = 1532.5
$ws.Range("N116").Value = -39124.668
$ws.Range("H122").Value = 5041.3335
$ws.Range("I122").Value = 5062
$ws.Range("K122").Value = 15186
$ws.Range("M122").Value = -12736
$ws.Range("H132").Value = 3265.8462
$ws.Range("I132").Value = 1934.3684
$ws.Range("J132").Value = 6879.857
$ws.Range("K132").Value = 5803.1052
$ws.Range("L132").Value = 20639.571
$ws.Range("M132").Value = -3273.1052
$ws.Range("N132").Value = -25699.571
$ws.Range("H136").Value = 7129.2964
$ws.Range("I136").Value = 6342.1055
$ws.Range("J136").Value = 8998.875
$ws.Range("K136").Value = 19026.3165
$ws.Range("L136").Value = 26996.625
$ws.Range("M136").Value = -16476.3165
$ws.Range("N136").Value = -32096.625
$ws.Range("H139").Value = 68425.25
$ws.Range("J139").Value = 68425.25
$ws.Range("L139").Value = 68425.25
$ws.Range("N139").Value = -78705.25
$ws.Range("M19").ClearContents()

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7516.533
$ws.Range("I3").Value = 761.5
$ws.Range("J3").Value = 34536.668
$ws.Range("K3").Value = 761.5
$ws.Range("L3").Value = 34536.668
$ws.Range("M3").Value = -647.5
$ws.Range("N3").Value = -34764.668
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("H134").Value = 2525.6667
$ws.Range("I134").Value = 1715.5
$ws.Range("K134").Value = 5146.5
$ws.Range("M134").Value = -2611.5
$ws.Range("M54").ClearContents()

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3292.818
$ws.Range("I58").Value = 1604
$ws.Range("J58").Value = 6248.25
$ws.Range("K58").Value = 1604
$ws.Range("L58").Value = 6248.25
$ws.Range("M58").Value = -1401
$ws.Range("N58").Value = -6654.25
$ws.Range("H94").Value = 2357.1667
$ws.Range("I94").Value = 1406.5
$ws.Range("K94").Value = 1406.5
$ws.Range("M94").Value = -955.5
$ws.Range("H132").Value = 3941.0938
$ws.Range("I132").Value = 3035.5
$ws.Range("K132").Value = 9106.5
$ws.Range("M132").Value = -6576.5
$ws.Range("H134").Value = 2701.4075
$ws.Range("I134").Value = 1580.1666
$ws.Range("J134").Value = 11671.333
$ws.Range("K134").Value = 4740.4998
$ws.Range("L134").Value = 35013.999
$ws.Range("M134").Value = -2205.4998
$ws.Range("N134").Value = -40083.999
$ws.Range("H135").Value = 63132
$ws.Range("J135").Value = 63132
$ws.Range("L135").Value = 63132
$ws.Range("N135").Value = -73272
$ws.Range("H136").Value = 3292.818
$ws.Range("I136").Value = 1604
$ws.Range("J136").Value = 6248.25
$ws.Range("K136").Value = 4812
$ws.Range("L136").Value = 18744.75
$ws.Range("M136").Value = -2262
$ws.Range("N136").Value = -23844.75
$ws.Range("H140").Value = 98589.625
$ws.Range("J140").Value = 98589.625
$ws.Range("L140").Value = 98589.625
$ws.Range("N140").Value = -108949.625

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2326.75
$ws.Range("J55").Value = 6800
$ws.Range("L55").Value = 20400
$ws.Range("N55").Value = -20754
$ws.Range("H126").Value = 3440
$ws.Range("I126").Value = 1150
$ws.Range("K126").Value = 3450
$ws.Range("M126").Value = 1490

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1440.2858
$ws.Range("I3").Value = 332.22223
$ws.Range("J3").Value = 3434.8
$ws.Range("K3").Value = 332.22223
$ws.Range("L3").Value = 3434.8
$ws.Range("M3").Value = -216.22223
$ws.Range("N3").Value = -3666.8
$ws.Range("H19").Value = 50
$ws.Range("I19").Value = 50
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 50
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 238
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("H58").Value = 45999.5
$ws.Range("J58").Value = 45999.5
$ws.Range("L58").Value = 45999.5
$ws.Range("N58").Value = -46553.5
$ws.Range("H62").Value = 44416.75
$ws.Range("I62").Value = 38834.5
$ws.Range("K62").Value = 38834.5
$ws.Range("M62").Value = -38148.5
$ws.Range("H65").Value = 44416.75
$ws.Range("I65").Value = 38834.5
$ws.Range("K65").Value = 116503.5
$ws.Range("M65").Value = -113071.5
$ws.Range("H102").Value = 2274.6216
$ws.Range("I102").Value = 1016
$ws.Range("K102").Value = 1016
$ws.Range("M102").Value = 606
$ws.Range("N19").ClearContents()
$ws.Range("N44").ClearContents()

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 3558.4
$ws.Range("I32").Value = 3558.4
$ws.Range("K32").Value = 3558.4
$ws.Range("M32").Value = -3241.4
$ws.Range("H55").Value = 1353217.2
$ws.Range("I55").Value = 2084303.9
$ws.Range("J55").Value = 3518.7693
$ws.Range("K55").Value = 2084303.9
$ws.Range("L55").Value = 3518.7693
$ws.Range("M55").Value = -2084130.9
$ws.Range("N55").Value = -3864.7693
$ws.Range("H61").Value = 5198.9585
$ws.Range("I61").Value = 4329.3687
$ws.Range("K61").Value = 4329.3687
$ws.Range("M61").Value = -4127.3687
$ws.Range("H113").Value = 5198.9585
$ws.Range("I113").Value = 4329.3687
$ws.Range("K113").Value = 4329.3687
$ws.Range("M113").Value = -2159.3687
$ws.Range("H132").Value = 4969.619
$ws.Range("I132").Value = 3446
$ws.Range("K132").Value = 10338
$ws.Range("M132").Value = -7808
$ws.Range("H136").Value = 4147.9
$ws.Range("I136").Value = 2286.6365
$ws.Range("J136").Value = 6422.778
$ws.Range("K136").Value = 6859.9095
$ws.Range("L136").Value = 19268.334
$ws.Range("M136").Value = -4309.9095
$ws.Range("N136").Value = -24368.334

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("H136").Value = 5141.1875
$ws.Range("I136").Value = 4053.2222
$ws.Range("J136").Value = 6540
$ws.Range("K136").Value = 12159.6666
$ws.Range("L136").Value = 19620
$ws.Range("M136").Value = -9609.6666
$ws.Range("N136").Value = -24720
$ws.Range("M12").ClearContents()
